$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 74.97871205021183
$ws.Range("C2").Value = 76.08898233642797
$ws.Range("D2").Value = 81.70912027889234
$ws.Range("E2").Value = 80.33317379989258
$ws.Range("B3").Value = 86.10921159065869
$ws.Range("C3").Value = 86.48485917783529
$ws.Range("D3").Value = 86.14888549446351
$ws.Range("E3").Value = 86.4564306229955
$ws.Range("B4").Value = 99.29849884191349
$ws.Range("C4").Value = 99.2754554553667
$ws.Range("D4").Value = 99.32840760905329
$ws.Range("E4").Value = 99.35672870082082
$ws.Range("B5").Value = 98.95797991079844
$ws.Range("C5").Value = 98.95965141652303
$ws.Range("D5").Value = 98.89240116085847
$ws.Range("E5").Value = 98.91056468915318
$ws.Range("B6").Value = 98.54715166863556
$ws.Range("C6").Value = 98.48310757927128
$ws.Range("D6").Value = 98.49366160461102
$ws.Range("E6").Value = 98.44194556324369
$ws.Range("B7").Value = 97.5406041668749
$ws.Range("C7").Value = 97.55052198775684
$ws.Range("D7").Value = 97.60325662890834
$ws.Range("E7").Value = 97.54219016892165
$ws.Range("B8").Value = 96.17749759138665
$ws.Range("C8").Value = 96.15619465059127
$ws.Range("D8").Value = 96.10886151195402
$ws.Range("E8").Value = 96.12392957128026
